# Footer block and Navvigation
#
# - C16 / C17: mark as text with an empty value (quote-prefix, "'" in the UI)
# - C51: "Message box Img and Department block"
# - C52: "Footer and Navigation"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C16").Value = "'"
$ws.Range("C17").Value = "'"

$ws.Range("C51").Value = "Message box Img and Department block"
$ws.Range("C52").Value = "Footer and Navigation"
